$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate data rows 2-8 (timestamps + Max/Avg/Min measurements)
$ws.Range("A2").Value = "20.08.2022 21:43:00"
$ws.Range("B2").Value = 9.759
$ws.Range("C2").Value = 5.325
$ws.Range("D2").Value = 0.786

$ws.Range("A3").Value = "20.08.2022 21:44:00"
$ws.Range("B3").Value = 9.403
$ws.Range("C3").Value = 3.960888
$ws.Range("D3").Value = 0.357

$ws.Range("A4").Value = "20.08.2022 21:45:00"
$ws.Range("B4").Value = 6.381
$ws.Range("C4").Value = 2.901444
$ws.Range("D4").Value = 0.526

$ws.Range("A5").Value = "20.08.2022 21:46:00"
$ws.Range("B5").Value = 8.899
$ws.Range("C5").Value = 4.654777
$ws.Range("D5").Value = 0.201

$ws.Range("A6").Value = "20.08.2022 21:47:00"
$ws.Range("B6").Value = 9.607
$ws.Range("C6").Value = 6.057666
$ws.Range("D6").Value = 1.733

$ws.Range("A7").Value = "20.08.2022 23:06:00"
$ws.Range("B7").Value = 7.506
$ws.Range("C7").Value = 4.092111
$ws.Range("D7").Value = 1.466

$ws.Range("A8").Value = "20.08.2022 23:15:00"
$ws.Range("B8").Value = 9.555
$ws.Range("C8").Value = 6.769777
$ws.Range("D8").Value = 1.271

# Reposition/resize the chart (graphicFrame) to its new anchor location
$co = $ws.ChartObjects(1)
$co.Left = 0
$co.Top = 120
$co.Width = 5012.1875
$co.Height = 750
